$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($range, $value)
    $c = $ws.Range($range)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextCell "D2" "37.597.28"
$ws.Range("E2").Value = "  -0.66%  "

# Row 3 - Ethereum
Set-TextCell "D3" "2.073.06"
$ws.Range("E3").Value = "  -0.53%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.06%  "

# Row 5 - BNB
Set-TextCell "D5" "231.77"
$ws.Range("E5").Value = "  -0.79%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -0.46%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.01%  "

# Row 8 - Solana
Set-TextCell "D8" "57.96"
$ws.Range("E8").Value = "  -1.94%  "

# Row 9 - Cardano
Set-TextCell "D9" "0.387"
$ws.Range("E9").Value = "  -2.13%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -1.77%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.36%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextCell "D12" "2.376.51"
$ws.Range("E12").Value = "  -0.66%  "

# Row 13 - Chainlink
Set-TextCell "D13" "14.75"
$ws.Range("E13").Value = "  -0.13%  "

# Row 14 - Avalanche
Set-TextCell "D14" "21.20"
$ws.Range("E14").Value = "  -0.45%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  -1.32%  "

# Row 16 - Polkadot
$ws.Range("E16").Value = "  -0.34%  "

# Row 17 - WrappedEther
Set-TextCell "D17" "2.077.33"
$ws.Range("E17").Value = "  -0.99%  "

# Row 18 - WrappedBTC
Set-TextCell "D18" "37.521.84"
$ws.Range("E18").Value = "  -0.67%  "

# Row 19 - Uniswap
$ws.Range("E19").Value = "  -0.12%  "

# Row 20 - Litecoin
Set-TextCell "D20" "70.01"
$ws.Range("E20").Value = "  -2.32%  "

# Row 21 - ShibaInu
Set-TextCell "D21" "0.0₃0827"
$ws.Range("E21").Value = "  -2.67%  "

# Row 22 - BitcoinCash
Set-TextCell "D22" "227.20"
$ws.Range("E22").Value = "  -0.42%  "

# Row 24 - Toncoin
Set-TextCell "D24" "2.40"
$ws.Range("E24").Value = "  +0.28%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  -3.08%  "

# Row 26 - Cosmos
Set-TextCell "D26" "9.88"
$ws.Range("E26").Value = "  +2.21%  "

# Row 27 - Monero
Set-TextCell "D27" "169.55"
$ws.Range("E27").Value = "  -0.97%  "

# Row 28 - Kaspa
$ws.Range("E28").Value = "  -5.01%  "

# Row 29 - EthereumClassic
Set-TextCell "D29" "19.36"
$ws.Range("E29").Value = "  -0.84%  "

# Row 30 - ImmutableX
$ws.Range("E30").Value = "  -4.48%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  -3.66%  "

# Row 33 - Hedera
$ws.Range("E33").Value = "  -1.34%  "

# Row 34 - InternetComputer(DFINITY)
Set-TextCell "D34" "4.63"
$ws.Range("E34").Value = "  -1.05%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  +1.08%  "

# Row 36 - WEMIXToken
$ws.Range("E36").Value = "  +0.64%  "

# Row 37 - RenderToken
$ws.Range("E37").Value = "  -4.09%  "

# Row 38 - BinanceUSD
$ws.Range("E38").Value = "  -0.03%  "

# Row 39 - THORChain
$ws.Range("E39").Value = "  -1.90%  "

# Row 40 - VeChain
$ws.Range("E40").Value = "  +3.82%  "

# Row 41 - Aave
Set-TextCell "D41" "98.13"
$ws.Range("E41").Value = "  -1.18%  "

# Row 42 - Cronos
Set-TextCell "D42" "0.0955"
$ws.Range("E42").Value = "  -3.05%  "

# Row 43 - Maker
Set-TextCell "D43" "1.486.58"
$ws.Range("E43").Value = "  +2.57%  "

# Row 44 - HuobiToken
Set-TextCell "D44" "2.91"
$ws.Range("E44").Value = "  +0.37%  "

# Row 45 - TrustWalletToken
$ws.Range("E45").Value = "  +2.65%  "

# Row 46 - InjectiveProtocol
Set-TextCell "D46" "16.56"
$ws.Range("E46").Value = "  -3.79%  "

# Row 47 - was ARBITRUM, now FTXToken (rows 47/48 swapped)
$ws.Range("B47").Value = "FTXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextCell "D47" "4.06"
$ws.Range("E47").Value = "  -2.14%  "

# Row 48 - was FTXToken, now ARBITRUM (rows 47/48 swapped)
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell "D48" "1.04"
$ws.Range("E48").Value = "  -2.53%  "

# Row 49 - FraxShare
$ws.Range("E49").Value = "  -1.83%  "

# Row 50 - MXToken
Set-TextCell "D50" "2.96"
$ws.Range("E50").Value = "  -1.13%  "

# Row 51 - RocketPoolETH
Set-TextCell "D51" "2.261.12"
$ws.Range("E51").Value = "  -0.72%  "
